$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell E1 - match the style of the existing header row (A1:D1)
$ws.Range("E1").Value = "Praat Label"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill "fear" as the Praat Label for every data row (2-31)
for ($i = 2; $i -le 31; $i++) {
    $ws.Cells.Item($i, 5).Value = "fear"
}
